$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Round row 5 (the last remaining data row after row 6 is deleted)
# values down from 3 decimal places of precision to 2 ---
for ($c = 2; $c -le 34; $c++) {
    $cell = $ws.Cells.Item(5, $c)
    $cell.Value = $excel.WorksheetFunction.Round($cell.Value2, 2)
}

# --- Narrow a handful of columns (Q=17, AA=27, AB=28, AC=29) from
# width 8 to width 7 (i.e. shrink the stored column width by one unit) ---
$narrowCols = @(17, 27, 28, 29)
foreach ($col in $narrowCols) {
    $ws.Columns.Item($col).ColumnWidth = $ws.Columns.Item($col).ColumnWidth - 1
}

# --- Remove the last data row (row 6); this also shrinks the sheet
# dimension from A1:AH6 down to A1:AH5 ---
$ws.Rows.Item(6).Delete()
